# Weekly update: rotate the price/quality figures for each row while
# leaving the market/product identification columns untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44230
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 160
$ws.Range("N2").Value = 16500
$ws.Range("O2").Value = 17000
$ws.Range("P2").Value = 16750
$ws.Range("Q2").Value = "`$/caja 18 kilos granel"
$ws.Range("S2").Value = 931
$ws.Range("T2").Value = 18

# Row 3
$ws.Range("D3").Value = 44230
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 160
$ws.Range("N3").Value = 14500
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 14750
$ws.Range("Q3").Value = "`$/caja 18 kilos granel"
$ws.Range("S3").Value = 819
$ws.Range("T3").Value = 18

# Row 4
$ws.Range("D4").Value = 44209
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 300
$ws.Range("N4").Value = 15500
$ws.Range("O4").Value = 16000
$ws.Range("P4").Value = 15750
$ws.Range("Q4").Value = "`$/caja 16 kilos granel"
$ws.Range("S4").Value = 984
$ws.Range("T4").Value = 16

# Row 5
$ws.Range("D5").Value = 44224
$ws.Range("L5").Value = "Especial"
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 16500
$ws.Range("O5").Value = 17000
$ws.Range("P5").Value = 16750
$ws.Range("Q5").Value = "`$/caja 16 kilos granel"
$ws.Range("S5").Value = 1047
$ws.Range("T5").Value = 16

# Row 6
$ws.Range("D6").Value = 44224
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 200
$ws.Range("N6").Value = 14500
$ws.Range("O6").Value = 15000
$ws.Range("P6").Value = 14750
$ws.Range("Q6").Value = "`$/caja 16 kilos granel"
$ws.Range("S6").Value = 922
$ws.Range("T6").Value = 16

# Row 7
$ws.Range("D7").Value = 44224
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 200
$ws.Range("N7").Value = 12500
$ws.Range("O7").Value = 13000
$ws.Range("P7").Value = 12750
$ws.Range("Q7").Value = "`$/caja 16 kilos granel"
$ws.Range("S7").Value = 797
$ws.Range("T7").Value = 16

# Row 8
$ws.Range("D8").Value = 44210
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 240
$ws.Range("N8").Value = 15500
$ws.Range("O8").Value = 16000
$ws.Range("P8").Value = 15750
$ws.Range("Q8").Value = "`$/caja 16 kilos granel"
$ws.Range("S8").Value = 984
$ws.Range("T8").Value = 16

# Row 9
$ws.Range("D9").Value = 44210
$ws.Range("L9").Value = "Segunda"
$ws.Range("M9").Value = 300
$ws.Range("N9").Value = 12500
$ws.Range("O9").Value = 13000
$ws.Range("P9").Value = 12750
$ws.Range("Q9").Value = "`$/caja 16 kilos granel"
$ws.Range("S9").Value = 797
$ws.Range("T9").Value = 16
